$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell used to stage new "numeric-looking" text values so that, when
# pasted (values-only) into the target cells, Excel keeps them as shared-string
# text cells (matching the original t="s" cell type) instead of silently
# re-typing them as numbers.
$helper = $ws.Range("AZ1000")
$helper.NumberFormat = "@"

# --- short-url column (B) : "P3pKuT" -> "7AkEyF" for every data row ---
$helper.Value = "7AkEyF"
$helper.Copy()
$lastRow = $ws.Range("A1").End(-4121).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).PasteSpecial(-4163)
}

# --- row 117 : refugees 8253 -> 9902, asylum_seekers 2115 -> 3029,
#     returned_refugees 32 -> 66 ---
$helper.Value = "9902"
$helper.Copy()
$ws.Range("N117").PasteSpecial(-4163)

$helper.Value = "3029"
$helper.Copy()
$ws.Range("O117").PasteSpecial(-4163)

$helper.Value = "66"
$helper.Copy()
$ws.Range("P117").PasteSpecial(-4163)

# --- row 120 : stateless 3629 -> 4466 ---
$helper.Value = "4466"
$helper.Copy()
$ws.Range("S120").PasteSpecial(-4163)

# --- row 121 : refugees 2775 -> 2250 ---
$helper.Value = "2250"
$helper.Copy()
$ws.Range("N121").PasteSpecial(-4163)

# Clean up the helper cell so it leaves no trace in the saved workbook.
$helper.ClearContents()
$helper.NumberFormat = "General"
$excel.CutCopyMode = $false
